# Update cryptos list with latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.014.88"
$ws.Range("E2").Value = '  -4.39%  '
$ws.Range("D3").Value = "'2.608.47"
$ws.Range("E3").Value = '  -3.47%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = "'517.24"
$ws.Range("E5").Value = '  -1.79%  '
$ws.Range("D6").Value = "'142.34"
$ws.Range("E6").Value = '  -2.26%  '
$ws.Range("E7").Value = '  +0.30%  '
$ws.Range("D8").Value = "'0.567"
$ws.Range("E8").Value = '  -1.82%  '
$ws.Range("D9").Value = "'6.70"
$ws.Range("E9").Value = '  -0.91%  '
$ws.Range("D10").Value = "'0.103"
$ws.Range("E10").Value = '  -2.93%  '
$ws.Range("E11").Value = '  -0.79%  '
$ws.Range("E12").Value = '  +1.03%  '
$ws.Range("D13").Value = "'3.067.79"
$ws.Range("E13").Value = '  -3.43%  '
$ws.Range("D14").Value = "'57.999.51"
$ws.Range("E14").Value = '  -4.35%  '
$ws.Range("D15").Value = "'20.91"
$ws.Range("E15").Value = '  -1.77%  '
$ws.Range("D17").Value = "'2.612.60"
$ws.Range("E17").Value = '  -3.89%  '
$ws.Range("D18").Value = "'4.40"
$ws.Range("E18").Value = '  -2.57%  '
$ws.Range("D19").Value = "'334.23"
$ws.Range("E19").Value = '  -3.40%  '
$ws.Range("D20").Value = "'10.33"
$ws.Range("E20").Value = '  -2.69%  '
$ws.Range("D21").Value = "'6.25"
$ws.Range("E21").Value = '  -3.38%  '
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("D23").Value = "'63.95"
$ws.Range("E23").Value = '  +0.76%  '
$ws.Range("E24").Value = '  -1.88%  '
$ws.Range("D25").Value = "'0.165"
$ws.Range("E25").Value = '  -2.82%  '
$ws.Range("D27").Value = "'7.09"
$ws.Range("E27").Value = '  -2.77%  '
$ws.Range("D28").Value = "'0.0₃0785"
$ws.Range("E28").Value = '  -4.51%  '
$ws.Range("D29").Value = "'6.60"
$ws.Range("E29").Value = '  -3.39%  '
$ws.Range("D31").Value = "'1.58"
$ws.Range("E31").Value = '  -1.42%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").Value = "'18.68"
$ws.Range("E32").Value = '  -1.95%  '
$ws.Range("B33").Value = 'Monero'
$ws.Range("C33").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D33").Value = "'150.14"
$ws.Range("E33").Value = '  +0.13%  '
$ws.Range("E34").Value = '  -4.12%  '
$ws.Range("E35").Value = '  -5.60%  '
$ws.Range("D36").Value = "'0.898"
$ws.Range("E36").Value = '  -4.20%  '
$ws.Range("E37").Value = '  -1.55%  '
$ws.Range("D38").Value = "'0.838"
$ws.Range("E38").Value = '  -4.21%  '
$ws.Range("E39").Value = '  -6.12%  '
$ws.Range("D40").Value = "'3.60"
$ws.Range("E40").Value = '  -1.81%  '
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = '  +0.40%  '
$ws.Range("D42").Value = "'0.597"
$ws.Range("E42").Value = '  -2.34%  '
$ws.Range("D43").Value = "'0.0963"
$ws.Range("E43").Value = '  -2.51%  '
$ws.Range("D44").Value = "'267.64"
$ws.Range("E44").Value = '  -5.44%  '
$ws.Range("D46").Value = "'19.10"
$ws.Range("E46").Value = '  -4.88%  '
$ws.Range("E47").Value = '  -1.24%  '
$ws.Range("D48").Value = "'2.029.01"
$ws.Range("E48").Value = '  -5.61%  '
$ws.Range("D49").Value = "'0.0228"
$ws.Range("E49").Value = '  -2.00%  '
$ws.Range("D50").Value = "'4.61"
$ws.Range("E50").Value = '  -6.17%  '
$ws.Range("D51").Value = "'18.20"
$ws.Range("E51").Value = '  -4.60%  '
